# UnivariateEDA.docx update ("Updated univariate EDA ho")
#
# The R-code echo blocks are simplified:
#   hist(~age,data=ars,main="",xlab="Age (yrs)")
#       -> hist(~age,data=ars,xlab="Age (yrs)")
#   hist(age~sex,data=ars,xlab="Age (yrs)",col="gray90")
#       -> hist(age~sex,data=ars,xlab="Age (yrs)")
#
# Each "source code" line is split across many single-style w:r runs (one
# run per syntax-highlighted token), so rather than doing a text Find &
# Replace (which would coalesce the runs spanning the match into one run
# and lose the original per-token rStyle formatting), we locate the exact
# substring to remove with Find.Execute (replace = wdReplaceNone) and then
# call Range.Delete() on the hit. That trims/removes only the runs that
# fall inside the deleted span and leaves untouched runs on either side
# exactly as they were.

$quote = [char]34

$d = $word.ActiveDocument

# 1) Remove `main="",` from the first hist(...) call.
$r1 = $d.Content
$needle1 = "main=" + $quote + $quote + ","
$found1 = $r1.Find.Execute($needle1, $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found1) {
    $r1.Delete()
}

# 2) Remove `,col="gray90"` from the hist(age~sex, ...) call (the first
#    "col=" occurrence -- the later barplot(...) call keeps its col= arg).
$r2 = $d.Content
$needle2 = "," + "col=" + $quote + "gray90" + $quote
$found2 = $r2.Find.Execute($needle2, $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found2) {
    $r2.Delete()
}

Write-Output "main= removed: $found1; col=gray90 removed: $found2"
